# Applies the "Add files via upload" edit:
#  - Adds a new multiple-choice question (about Python's range()) to the
#    "3_" worksheet, which was previously empty.
#  - Updates the stray range selection left on the "4_MultC" worksheet
#    (A1:C6 -> A1:C5).
#  - Moves the active tab / tabSelected flag from "2_" to "3_" (the newly
#    filled-in sheet becomes the one the workbook opens on).

$wb = $excel.ActiveWorkbook

# --- 1. Tidy up the leftover selection on "4_MultC" (A1:C6 -> A1:C5). ----
$wsMultC = $wb.Worksheets.Item("4_MultC")
$wsMultC.Range("A1:C5").Select()

# --- 2. Fill in the new question table on sheet "3_". --------------------
$ws = $wb.Worksheets.Item("3_")

$ws.Range("A1").Value = "If a line of code reads 'for x in range(10,15), what are the smallest and largest values of 'x' in the loop?"
$ws.Range("B1").Value = "Correct"
$ws.Range("C1").Value = "Comment"

$ws.Range("A2").Value = "Smallest: 0; Largest: 15"
$ws.Range("B2").Value = "N"

$ws.Range("A3").Value = "Smallest: 11; Largest: 15"
$ws.Range("B3").Value = "N"

$ws.Range("A4").Value = "Smallest: 10; Largest: 14"
$ws.Range("B4").Value = "Y"
$ws.Range("C4").Value = "Yep!  The range starts at the first number, and stops one value before the second number."

$ws.Range("A5").Value = "Smallest: 10; Largest: 15"
$ws.Range("B5").Value = "N"

$ws.Range("A6").Value = "Smallest: 0; Largest: 14"
$ws.Range("B6").Value = "N"

# Row heights: header row (wrapped, tallest) and the correct-answer row
# (which also carries a wrapped comment) grow to fit their text.
$ws.Rows.Item(1).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 45

# --- 3. Make "3_" the active sheet/tab, with C5 selected. ----------------
$ws.Activate()
$ws.Range("C5").Select()
